# North Star POD page done
# Adds a new "North_star" worksheet at the end of the workbook, populates
# its summary table, and updates the previously-active sheet's selection.

$wb = $excel.ActiveWorkbook

# --- Update the selection on "Thought_leadership" (was the active sheet) ---
$wsThought = $wb.Worksheets.Item("Thought_leadership")
$wsThought.Range("C16:L20").Select()

# --- Add the new "North_star" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "North_star"

# Header row (row 14)
$ws.Range("E14").Value = "id"
$ws.Range("G14").Value = "label"
$ws.Range("I14").Value = "percentage"
$ws.Range("K14").Value = "title"
$ws.Range("N14").Value = "heading"

# Data row 16
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = "Value Per Quarter"
$ws.Range("I16").Value = 76
$ws.Range("K16").Value = "EN $"
$ws.Range("N16").Value = "Total Revenue: 2023-2024"

# Data row 17
$ws.Range("E17").Value = 2
$ws.Range("G17").Value = "Measured in $ Value"
$ws.Range("I17").Value = 76
$ws.Range("K17").Value = "Customer Value ADD - "
$ws.Range("N17").Value = "Total Revenue: 2023-2024"

# Data row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = "Attrition %"
$ws.Range("I18").Value = 60
$ws.Range("N18").Value = "Total Revenue: 2023-2024"

# Make North_star the active tab/selection, matching the saved view state
$ws.Activate()
$ws.Range("Q17").Select()
